$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = 2471
$ws.Range("B59").Value = 7020
$ws.Range("C59").Value = "BENEFÍCIOS PREVIDENCIÁRIOS - ARTEMIG"
